$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2..407) holds the "Förändrad" (changed) date as an Excel
# serial number. All of them were bumped by one day: 45180 -> 45181.
$ws.Range("C2:C407").Value = 45181
